$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 242 - this shifts the existing rows
# 242-300 down to 243-301 (dimension grows from R300 to R301).
$ws.Rows("242:242").Insert()

# Populate the newly inserted row 242 with its data. Columns A, B, C, E,
# F, G, H, I, N, Q, R keep the same values the (old) row 242 had; D, J,
# K, L, M, O, P take the new values from the diff.
$ws.Cells.Item(242, 1).Value = 11
$ws.Cells.Item(242, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(242, 3).Value = "Bíobío"
$ws.Cells.Item(242, 4).Value = 44889
$ws.Cells.Item(242, 5).Value = 8
$ws.Cells.Item(242, 6).Value = 100114013
$ws.Cells.Item(242, 7).Value = "Zanahoria"
$ws.Cells.Item(242, 8).Value = "Sin especificar"
$ws.Cells.Item(242, 9).Value = "Primera"
$ws.Cells.Item(242, 10).Value = 180
$ws.Cells.Item(242, 11).Value = 7500
$ws.Cells.Item(242, 12).Value = 8000
$ws.Cells.Item(242, 13).Value = 7778
$ws.Cells.Item(242, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(242, 15).Value = "Región Metropolitana"
$ws.Cells.Item(242, 16).Value = 389
$ws.Cells.Item(242, 17).Value = 20
$ws.Cells.Item(242, 18).Value = "Hortaliza"
